$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover Aurora, a beautiful online slot game with 5 reels and 30 paylines. Play for free and find out about its gameplay, graphics, and bonuses.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2. Drop the duplicated bold "Play Aurora Slot..." paragraph that
#    used to sit near the end of the document (just before the
#    italic summary paragraph).
# ------------------------------------------------------------------
$oldTitleText = "Play Aurora Slot for Free - Review of Aurora Online Slot"
$dupPara = $null
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    $candText = $cand.Range.Text.TrimEnd([char]13)
    if ($candText -eq $oldTitleText) {
        $dupPara = $cand
    }
}
if ($dupPara -ne $null) {
    $dupPara.Range.Delete()
}

# ------------------------------------------------------------------
# 3. Replace the text of the last (italic) paragraph with the new
#    image-generation prompt, keeping its italic run formatting and
#    the paragraph's leading empty run untouched.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$oldSummaryText = "Discover Aurora, a beautiful online slot game with 5 reels and 30 paylines. Play for free and find out about its gameplay, graphics, and bonuses."
$promptText = 'Prompt: Please create a feature image that fits the "Aurora" online slot game. The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. DALLE, can you draw a feature image for the "Aurora" online slot game? The image should be in a cartoon style and feature a happy Maya warrior wearing glasses. The warrior should be standing in front of a vibrant backdrop of the northern lights, with the symbols of the game appearing above or around them. The image should be playful and eye-catching, with bright colors and energetic lines. The Maya warrior should look like they''re enjoying playing the game and celebrating a win. Please make sure the image fits the theme of the game and will appeal to both experienced and novice players. Thank you!'

$lastParaText = $lastPara.Range.Text.TrimEnd([char]13)
if ($lastParaText -eq $oldSummaryText) {
    $target = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
} else {
    $target = $lastPara.Range
}

$promptEscaped = $promptText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;").Replace('"', "&quot;")
$promptXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>$promptEscaped</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$target.InsertXML($promptXml)
